# Update countries & provincias Spain
# Applies the weekly COVID data refresh: updates the "last updated" timestamp,
# refreshes several countries' statistics, and re-labels rows whose country
# order shifted in the source data (which - because the worksheet stores one
# country name per row sequentially - shows up as swapped country names
# between adjacent rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 00:52"

# --- Row 4: Estados Unidos (stats refresh only) ----------------------------
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 394587
$ws.Range("C4").Value = 27583
$ws.Range("D4").Value = 21674
$ws.Range("E4").Value = 360165
$ws.Range("F4").Value = 9169
$ws.Range("G4").Value = 1877
$ws.Range("H4").Value = 12748

# --- Row 30: Chequia (stats refresh only) ----------------------------------
$ws.Range("A30").Value = "Chequia"
$ws.Range("B30").Value = 5017
$ws.Range("C30").Value = 195
$ws.Range("D30").Value = 172
$ws.Range("E30").Value = 4757
$ws.Range("F30").Value = 86
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 88

# --- Rows 95/96: Jordania & Reunion swap places, stats refresh -------------
$ws.Range("A95").Value = "Reunion"
$ws.Range("B95").Value = 358
$ws.Range("C95").Value = 9
$ws.Range("D95").Value = 40
$ws.Range("E95").Value = 318
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0

$ws.Range("A96").Value = "Jordania"
$ws.Range("B96").Value = 353
$ws.Range("C96").Value = 4
$ws.Range("D96").Value = 138
$ws.Range("E96").Value = 209
$ws.Range("F96").Value = 5
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 6

# --- Rows 116/117: Venezuela & Mayotte swap places, stats refresh ----------
$ws.Range("A116").Value = "Mayotte"
$ws.Range("B116").Value = 171
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 22
$ws.Range("E116").Value = 147
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 2

$ws.Range("A117").Value = "Venezuela"
$ws.Range("B117").Value = 165
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 65
$ws.Range("E117").Value = 93
$ws.Range("F117").Value = 6
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

# --- Rows 193-196: San Vicente y las Granadinas moves up in front of ------
# --- Malaui, Santa Sede and Belice shift down one row, stats refresh ------
$ws.Range("A193").Value = "San Vicente y las Granadinas"
$ws.Range("B193").Value = 8
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 7
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Malaui"
$ws.Range("B194").Value = 8
$ws.Range("C194").Value = 3
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 1
$ws.Range("G194").Value = 1
$ws.Range("H194").Value = 1

$ws.Range("A195").Value = "Santa Sede"
$ws.Range("B195").Value = 7
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 7
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Belice"
$ws.Range("B196").Value = 7
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 6
$ws.Range("F196").Value = 1
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1
